# Adds a new subscriber row (row 33) to the subscriber list sheet and
# drops the stray empty E32 cell, matching the "adding combobox for
# selecting the sectoral unions" commit: a new entry whose "region"
# column picks the option "4RE" and whose "specialty" column is the
# Arabic value "الماشية".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C and D hold long digit strings ("12345678") that must stay
# text (not be auto-coerced to numbers) to match the rest of the sheet.
$ws.Range("C33:D33").NumberFormat = "@"

$ws.Range("A33").Value = "ddd"
$ws.Range("B33").Value = "zdz"
$ws.Range("C33").Value = "12345678"
$ws.Range("D33").Value = "12345678"
$ws.Range("E33").Value = "4RE"
$ws.Range("F33").Value = "الماشية"

# Writing the new row causes the engine to materialize an empty E32
# placeholder as the used range grows; clear it so the cell disappears
# again, same as in the target sheet.
$ws.Range("E32").ClearContents()
